# Add a new row to the "Completed" reading list sheet for
# "The Case Against Socialism" by Rand Paul.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$newRow = 82

$ws.Cells.Item($newRow, 1).Value = "The Case Against Socialism"
$ws.Cells.Item($newRow, 2).Value = "Rand Paul"

# Start Date / Finish Date: write the serial values then copy the
# existing date-formatted cells' number format down so the new cells
# reuse the workbook's existing date style instead of creating a new one.
$ws.Cells.Item($newRow, 3).Value = 43980
$ws.Cells.Item($newRow, 4).Value = 43983

$ws.Range("C81").Copy() | Out-Null
$ws.Range("C82").PasteSpecial(-4122) | Out-Null
$ws.Range("D81").Copy() | Out-Null
$ws.Range("D82").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 5).Value = "socialism;politics;economics;capitalism"
$ws.Cells.Item($newRow, 6).Value = "Audio"
$ws.Cells.Item($newRow, 7).Value = "10 Hours 19 Mins"

$ws.Range("A83").Select()
